# Logica di append per superdettagli
# Replace the placeholder "next id" values in column AA with the real
# sequential super-detail ids (1..12), and move the selection to the
# freshly-edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $ws.Range("AA$row").Value = $values[$i]
}

$ws.Range("AA4:AA15").Select()
